$d = $word.ActiveDocument

# Update the date line (first paragraph)
$d.Content.Find.Execute("2024-06-26 Wednesday", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "2024-06-27 Thursday", 2) | Out-Null

# Update the 5x5 grid of division problems. The table has 20 rows; content
# lives in rows 1, 5, 9, 13, 17 (5 columns each), the rest are spacer rows.
$t = $d.Tables.Item(1)

$values = @(
    @("29÷8=3, 5",  "48÷5=9, 3",  "64÷2=32, 0", "13÷2=6, 1",  "21÷7=3, 0"),
    @("11÷6=1, 5",  "92÷2=46, 0", "64÷2=32, 0", "15÷9=1, 6",  "23÷6=3, 5"),
    @("84÷8=10, 4", "79÷9=8, 7",  "10÷4=2, 2",  "37÷5=7, 2",  "43÷9=4, 7"),
    @("41÷5=8, 1",  "16÷9=1, 7",  "58÷4=14, 2", "66÷6=11, 0", "65÷3=21, 2"),
    @("91÷2=45, 1", "51÷3=17, 0", "31÷3=10, 1", "89÷7=12, 5", "75÷9=8, 3")
)

$rowIndexes = @(1, 5, 9, 13, 17)

for ($i = 0; $i -lt 5; $i++) {
    $rowIdx = $rowIndexes[$i]
    for ($col = 1; $col -le 5; $col++) {
        $cell = $t.Cell($rowIdx, $col)
        $cell.Range.Text = $values[$i][$col - 1]
    }
}
